$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("H7").Value = "t=0.31s pour une charge du banc de condensateur d'environs 1700V"
$ws.Range("F7").Value = "t=0.31s, Pmoy = 2.53MW et Pmax=3.55MW pour PSIM, Pmoy=2.58MW et Pmax=3.56MW pour SPS"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
